$d = $word.ActiveDocument

function Replace-ExactText([string]$searchText, [string]$newText) {
    # Locate the exact text, then force a clean re-write of that range so that
    # any proofErr wrappers / run splits inside the found range collapse into
    # a single fresh run (mirrors what Word does when you select text and
    # retype it).
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.Text = $searchText
    $found = $rng.Find.Execute()
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return $null
    }
    # Two-step write: first collapse the range to a placeholder, then set the
    # real text. This avoids a silent no-op when old/new text happen to share
    # characters, and ensures the whole matched span re-merges into one run.
    $rng.Text = "@@TMP@@"
    $rng2 = $d.Range($rng.Start, $rng.End)
    $rng2.Text = $newText
    return $rng2
}

# 1) Remove the stray _GoBack bookmark that originally sat right after
#    "австрийского" (it gets relocated further down by this edit).
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
}

# 2) Fix the typo "факрам" -> "фактам" (drop the now-stale spelling proofErr
#    markers by rewriting the whole enclosing phrase).
Replace-ExactText "высказанным факрам." "высказанным фактам." | Out-Null

# 3) Fix the grammar slip "из содержание" -> "их содержание" (drop the
#    now-stale grammar proofErr markers the same way).
Replace-ExactText "определяет из содержание." "определяет их содержание." | Out-Null

# 4) The run that only contains the inline canvas drawing loses its explicit
#    ru-RU language override (it still keeps eastAsia="en-GB").
$shp = $d.InlineShapes.Item(1)
$shpPara = $shp.Range.Paragraphs.Item(1).Range
$shpPara.Font.LanguageID = "en-GB"

# 5) Merge the "Пример" / ": " runs into a single "Пример: " run.
Replace-ExactText "Пример: " "Пример: " | Out-Null

# 6) Fix "принцима атомизме" -> "принципа атомизма" (drop stale spelling
#    proofErr markers), then drop the new _GoBack bookmark right after the
#    fixed word, before the following comma.
Replace-ExactText "Развитие принцима атомизме," "Развитие принципа атомизма," | Out-Null

$rngFinal = $d.Content
$rngFinal.Find.ClearFormatting()
$rngFinal.Find.Forward = $true
$rngFinal.Find.Wrap = 0
$rngFinal.Find.MatchCase = $true
$rngFinal.Find.Text = "принципа атомизма"
$foundFinal = $rngFinal.Find.Execute()
if ($foundFinal) {
    $bmRng = $d.Range($rngFinal.End, $rngFinal.End)
    $d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
}
